# Redefine the base/default RDF prefix from the empty/colon-only name (":")
# to the explicit name "ome" throughout the workbook.
#
# This touches:
#   1. The "@prefix" sheet, where row 1 declared the base prefix with a
#      blank name ("" / ":") -> now named "ome".
#   2. Every other sheet ("Image", "Pixels", "Channel", "Color",
#      "Binary_Data") that referenced that base prefix via a leading
#      colon (":pixels", ":Image", ":Pixels", ... ) -> now "ome:pixels",
#      "ome:Image", "ome:Pixels", ...

$wb = $excel.ActiveWorkbook

# ---- @prefix sheet: name the base prefix "ome" ----
$wsPrefix = $wb.Worksheets.Item("@prefix")
$wsPrefix.Range("A1").Value = "ome"

# ---- Image sheet ----
$wsImage = $wb.Worksheets.Item("Image")
$wsImage.Range("E3").Value = "ome:pixels"
$wsImage.Range("F3").Value = "ome:acquisitionDate"
$wsImage.Range("B4").Value = "ome:Image"
$wsImage.Range("E4").Value = "ome:Pixels"

# ---- Pixels sheet ----
$wsPixels = $wb.Worksheets.Item("Pixels")
$wsPixels.Range("D3").Value = "ome:pixelType"
$wsPixels.Range("E3").Value = "ome:dimensionOrder"
$wsPixels.Range("F3").Value = "ome:physicalSizeX"
$wsPixels.Range("G3").Value = "ome:physicalSizeY"
$wsPixels.Range("H3").Value = "ome:sizeC"
$wsPixels.Range("I3").Value = "ome:sizeT"
$wsPixels.Range("J3").Value = "ome:sizeX"
$wsPixels.Range("K3").Value = "ome:sizeY"
$wsPixels.Range("L3").Value = "ome:sizeZ"
$wsPixels.Range("M3").Value = "ome:channel"
$wsPixels.Range("N3").Value = "ome:binData"
$wsPixels.Range("B4").Value = "ome:Pixels"
$wsPixels.Range("D4").Value = "ome:PixelType"
$wsPixels.Range("E4").Value = "ome:DimensionOrder"
$wsPixels.Range("M4").Value = "ome:Channel"
$wsPixels.Range("N4").Value = "ome:BinData"

# ---- Channel sheet ----
$wsChannel = $wb.Worksheets.Item("Channel")
$wsChannel.Range("D3").Value = "ome:color"
$wsChannel.Range("B4").Value = "ome:Channel"
$wsChannel.Range("D4").Value = "ome:Color"

# ---- Color sheet ----
$wsColor = $wb.Worksheets.Item("Color")
$wsColor.Range("B4").Value = "ome:Color"

# ---- Binary_Data sheet ----
$wsBinData = $wb.Worksheets.Item("Binary_Data")
$wsBinData.Range("C3").Value = "ome:bigEndian"
$wsBinData.Range("D3").Value = "ome:data"
$wsBinData.Range("E3").Value = "ome:length"
$wsBinData.Range("B4").Value = "ome:BinData"
